# This workbook contains a single sheet of weekly market-price records for
# "Berenjena" (rows 2-138, one record per row). The update shifts every
# existing record (rows 48-138) down by one row to make room for a brand
# new record inserted at row 48, and the record that used to be the very
# last one (old row 138) becomes the new last row (139).
#
# Concretely:
#   new_row[48] = brand-new record
#   new_row[r]  = old_row[r-1]   for r = 49..139
#
# We replicate this with a single bulk block-copy (read A48:R138, write it
# one row lower into A49:R139) and then overwrite row 48 with the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the destination block for column D (the date column) already
# carries the same date/time number format as the rest of the column
# before we drop values into it, so Excel doesn't invent a new ad-hoc
# number format when it sees date-like serials.
$ws.Range("D49:D139").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Read the current block (old rows 48-138) then write it back one row
# down (new rows 49-139). This performs the down-shift for every column.
$srcBlock = $ws.Range("A48:R138").Value()
$ws.Range("A49:R139").Value = $srcBlock

# Now overwrite row 48 with the brand-new record. Only the fields that
# actually change are listed; the rest (A,B,C,E,F,G,H,I,N,Q,R) keep the
# values they already had (they're identical for every row in this
# sheet, so row 48 still holds the correct data for them).
$ws.Range("D48").Value = 44973
$ws.Range("J48").Value = 170
$ws.Range("K48").Value = 8500
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = 8765
$ws.Range("O48").Value = "Región Metropolitana"
$ws.Range("P48").Value = 146
